# UZOVI en productgroep toegevoegd aan BGGZ trajecten datastructuur
#
# The BGGZ_trajecten list (numId 1003) currently reads:
#   dossierid (key), bggzid (key), startdatum, einddatum, status,
#   ZVZ_initieel, ZVZ_actueel, afsluitreden (...)
#
# Target list:
#   dossierid (key), bggzid (key), bsn (key), startdatum, einddatum,
#   trajectstatus, ZVZ_initieel, ZVZ_actueel, productgroep,
#   uzovi (...), afsluitreden (...)

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Insert "bsn (key)" right after "bggzid (key)" / before "startdatum"
#    (only inside the BGGZ_trajecten list - the other "bggzid (key)"
#    item, in BGGZ dossiers, is followed by "status", not "startdatum").
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$bggzidIdx = -1
for ($i = 1; $i -lt $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim() -eq "bggzid (key)") {
        $next = $paras.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "startdatum") {
            $bggzidIdx = $i
        }
    }
}
if ($bggzidIdx -eq -1) {
    throw "Could not locate 'bggzid (key)' -> 'startdatum' anchor in BGGZ_trajecten"
}

$startdatumPara = $d.Paragraphs.Item($bggzidIdx + 1)
$startdatumPara.Range.InsertParagraphBefore()
$bsnPara = $d.Paragraphs.Item($bggzidIdx + 1)
$bsnPara.Range.Text = "bsn (key)"

# ---------------------------------------------------------------
# 2. Rename "status" -> "trajectstatus" inside BGGZ_trajecten
#    (the "status" paragraph immediately followed by "ZVZ_initieel").
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$statusIdx = -1
for ($i = 1; $i -lt $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim() -eq "status") {
        $next = $paras.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "ZVZ_initieel") {
            $statusIdx = $i
        }
    }
}
if ($statusIdx -eq -1) {
    throw "Could not locate 'status' -> 'ZVZ_initieel' anchor in BGGZ_trajecten"
}

$statusPara = $d.Paragraphs.Item($statusIdx)
$statusPara.Range.Text = "trajectstatus"

# ---------------------------------------------------------------
# 3. Insert "productgroep" and the "uzovi (...)" item right after
#    "ZVZ_actueel" / before "afsluitreden (niet in queries...)".
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$zvzIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim() -eq "ZVZ_actueel") {
        $zvzIdx = $i
    }
}
if ($zvzIdx -eq -1) {
    throw "Could not locate 'ZVZ_actueel' anchor in BGGZ_trajecten"
}

$zvzPara = $d.Paragraphs.Item($zvzIdx)
$zvzPara.Range.InsertParagraphAfter()
$productgroepPara = $d.Paragraphs.Item($zvzIdx + 1)
$productgroepPara.Range.Text = "productgroep"

$productgroepPara = $d.Paragraphs.Item($zvzIdx + 1)
$productgroepPara.Range.InsertParagraphAfter()
$uzoviPara = $d.Paragraphs.Item($zvzIdx + 2)
$uzoviPara.Range.Text = "uzovi (moet eigenlijk van het dossier komen, maar geen goed bronbestand nu)"

Write-Output "done"
